$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 17143.2
$ws.Range("I21").Value = 11017
$ws.Range("J21").Value = 18674.75
$ws.Range("K21").Value = 11017
$ws.Range("L21").Value = 18674.75
$ws.Range("M21").Value = -10549
$ws.Range("N21").Value = -19610.75
$ws.Range("H23").Value = 17143.2
$ws.Range("I23").Value = 11017
$ws.Range("J23").Value = 18674.75
$ws.Range("K23").Value = 11017
$ws.Range("L23").Value = 18674.75
$ws.Range("M23").Value = -10783
$ws.Range("N23").Value = -19142.75
$ws.Range("H29").Value = 184.5
$ws.Range("I29").Value = 184.5
$ws.Range("K29").Value = 553.5
$ws.Range("M29").Value = -272.5
$ws.Range("H38").Value = 557.8125
$ws.Range("I38").Value = 443.75
$ws.Range("J38").Value = 900
$ws.Range("K38").Value = 1331.25
$ws.Range("L38").Value = 2700
$ws.Range("M38").Value = -959.25
$ws.Range("N38").Value = -3444
$ws.Range("H58").Value = 1179.4762
$ws.Range("J58").Value = 1983.3334
$ws.Range("L58").Value = 5950.0002
$ws.Range("N58").Value = -6250.0002
$ws.Range("H87").Value = 26363.637
$ws.Range("J87").Value = 26363.637
$ws.Range("L87").Value = 26363.637
$ws.Range("N87").Value = -28859.637
$ws.Range("H90").Value = 26363.637
$ws.Range("J90").Value = 26363.637
$ws.Range("L90").Value = 79090.91099999999
$ws.Range("N90").Value = -91570.91099999999
$ws.Range("H123").Value = 30000
$ws.Range("J123").Value = 30000
$ws.Range("L123").Value = 30000
$ws.Range("N123").Value = -39800
$ws.Range("H138").Value = 2095.093
$ws.Range("I138").Value = 2100.9
$ws.Range("J138").Value = 2090.0435
$ws.Range("K138").Value = 6302.700000000001
$ws.Range("L138").Value = 6270.130500000001
$ws.Range("M138").Value = -1162.700000000001
$ws.Range("N138").Value = -16550.1305

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 13028.111
$ws.Range("J23").Value = 9676.471
$ws.Range("L23").Value = 9676.471
$ws.Range("N23").Value = -10194.471
$ws.Range("H37").Value = 11992.6
$ws.Range("J37").Value = 11992.6
$ws.Range("L37").Value = 11992.6
$ws.Range("N37").Value = -12538.6
$ws.Range("H44").Value = 16398
$ws.Range("J44").Value = 16398
$ws.Range("L44").Value = 16398
$ws.Range("N44").Value = -17374
$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()
$ws.Range("H63").Value = 4558.2354
$ws.Range("I63").Value = 3842.5715
$ws.Range("J63").Value = 5059.2
$ws.Range("K63").Value = 3842.5715
$ws.Range("L63").Value = 5059.2
$ws.Range("M63").Value = -3156.5715
$ws.Range("N63").Value = -6431.2
$ws.Range("H66").Value = 4558.2354
$ws.Range("I66").Value = 3842.5715
$ws.Range("J66").Value = 5059.2
$ws.Range("K66").Value = 19212.8575
$ws.Range("L66").Value = 25296
$ws.Range("M66").Value = -15780.8575
$ws.Range("N66").Value = -32160
$ws.Range("H80").Value = 21665.666
$ws.Range("J80").Value = 21665.666
$ws.Range("L80").Value = 21665.666
$ws.Range("N80").Value = -23661.666
$ws.Range("H83").Value = 21665.666
$ws.Range("J83").Value = 21665.666
$ws.Range("L83").Value = 64996.99800000001
$ws.Range("N83").Value = -74980.99800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 12973.8
$ws.Range("I82").Value = 5576.75
$ws.Range("J82").Value = 21427.572
$ws.Range("K82").Value = 5576.75
$ws.Range("L82").Value = 21427.572
$ws.Range("M82").Value = -5193.75
$ws.Range("N82").Value = -22193.572
$ws.Range("H85").Value = 12973.8
$ws.Range("I85").Value = 5576.75
$ws.Range("J85").Value = 21427.572
$ws.Range("K85").Value = 5576.75
$ws.Range("L85").Value = 21427.572
$ws.Range("M85").Value = -4250.75
$ws.Range("N85").Value = -24079.572
$ws.Range("H107").Value = 1509.125
$ws.Range("I107").Value = 1373.1177
$ws.Range("J107").Value = 1839.4286
$ws.Range("K107").Value = 1373.1177
$ws.Range("L107").Value = 1839.4286
$ws.Range("M107").Value = 546.8823
$ws.Range("N107").Value = -5679.4286

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 29889.111
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 29889.111
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 29889.111
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -30113.111
$ws.Range("H122").Value = 1479.7091
$ws.Range("I122").Value = 1497.5714
$ws.Range("J122").Value = 1461.1852
$ws.Range("K122").Value = 4492.7142
$ws.Range("L122").Value = 4383.5556
$ws.Range("M122").Value = -2042.7142
$ws.Range("N122").Value = -9283.5556
$ws.Range("H132").Value = 6174787
$ws.Range("I132").Value = 1646.7727
$ws.Range("J132").Value = 33336604
$ws.Range("K132").Value = 4940.3181
$ws.Range("L132").Value = 100009812
$ws.Range("M132").Value = -2410.3181
$ws.Range("N132").Value = -100014872

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 10638749
$ws.Range("I34").Value = 165.875
$ws.Range("J34").Value = 12821022
$ws.Range("K34").Value = 497.625
$ws.Range("L34").Value = 38463066
$ws.Range("M34").Value = -413.625
$ws.Range("N34").Value = -38463234
$ws.Range("H39").Value = 2049.8
$ws.Range("J39").Value = 2049.8
$ws.Range("L39").Value = 6149.400000000001
$ws.Range("N39").Value = -6737.400000000001
$ws.Range("H55").Value = 1881.8182
$ws.Range("I55").Value = 1800
$ws.Range("J55").Value = 1890
$ws.Range("K55").Value = 5400
$ws.Range("L55").Value = 5670
$ws.Range("M55").Value = -5223
$ws.Range("N55").Value = -6024
$ws.Range("H122").Value = 6163.1665
$ws.Range("I122").Value = 545.6429000000001
$ws.Range("K122").Value = 4910.7861
$ws.Range("M122").Value = -2460.7861

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 4938.143
$ws.Range("J43").Value = 18166
$ws.Range("L43").Value = 18166
$ws.Range("N43").Value = -18468
$ws.Range("H46").Value = 4349.95
$ws.Range("J46").Value = 4349.95
$ws.Range("L46").Value = 4349.95
$ws.Range("N46").Value = -4661.95
$ws.Range("H57").Value = 10012666
$ws.Range("J57").Value = 10012666
$ws.Range("L57").Value = 10012666
$ws.Range("N57").Value = -10014306
$ws.Range("H80").Value = 63677880
$ws.Range("J80").Value = 201500
$ws.Range("L80").Value = 201500
$ws.Range("N80").Value = -203496
$ws.Range("H83").Value = 63677880
$ws.Range("J83").Value = 201500
$ws.Range("L83").Value = 1007500
$ws.Range("N83").Value = -1017484
$ws.Range("H132").Value = 3658.5264
$ws.Range("I132").Value = 3465.2856
$ws.Range("J132").Value = 4199.6
$ws.Range("K132").Value = 10395.8568
$ws.Range("L132").Value = 12598.8
$ws.Range("M132").Value = -7865.856800000001
$ws.Range("N132").Value = -17658.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 2550001
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 2550001
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 2550001
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -2550225
$ws.Range("H16").Value = 14287428
$ws.Range("I16").Value = 1162.8572
$ws.Range("K16").Value = 1162.8572
$ws.Range("M16").Value = -992.8571999999999
$ws.Range("H132").Value = 3572.8333
$ws.Range("I132").Value = 2879.7856
$ws.Range("K132").Value = 8639.356800000001
$ws.Range("M132").Value = -6109.356800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 70029
$ws.Range("J34").Value = 70029
$ws.Range("L34").Value = 70029
$ws.Range("N34").Value = -70435
$ws.Range("H42").Value = 70049
$ws.Range("J42").Value = 70049
$ws.Range("L42").Value = 70049
$ws.Range("N42").Value = -70805
$ws.Range("H43").Value = 30000
$ws.Range("J43").Value = 30000
$ws.Range("L43").Value = 30000
$ws.Range("N43").Value = -30298
$ws.Range("H81").Value = 6385.4443
$ws.Range("I81").Value = 5578.1665
$ws.Range("J81").Value = 8000
$ws.Range("K81").Value = 11156.333
$ws.Range("L81").Value = 16000
$ws.Range("M81").Value = -10095.333
$ws.Range("N81").Value = -18122
$ws.Range("H84").Value = 6385.4443
$ws.Range("I84").Value = 5578.1665
$ws.Range("J84").Value = 8000
$ws.Range("K84").Value = 55781.665
$ws.Range("L84").Value = 80000
$ws.Range("M84").Value = -50477.665
$ws.Range("N84").Value = -90608
$ws.Range("H100").Value = 2886.7144
$ws.Range("I100").Value = 3041.4
$ws.Range("J100").Value = 2500
$ws.Range("K100").Value = 6082.8
$ws.Range("L100").Value = 5000
$ws.Range("M100").Value = -5541.8
$ws.Range("N100").Value = -6082
